# Devauspaivakirja.xlsx - add a "Nav komponentti" diary entry (row 9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 9 by copying the formatting of row 2 (A2 carries the date
# number format used by the whole pvm column) so the new date cell reuses
# the existing cell style instead of registering a new one.
$ws.Range("A2").Copy($ws.Range("A9"))

# New diary entry: 16.3.2021, 100 minutes, "Nav komponentti".
$ws.Range("A9").Value = 44271
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = "Nav komponentti"

# Move/restore the cell selection like it was left after entering the row.
$ws.Range("C11").Select()
